# MasterStoryList.xlsx edit: add Notes/Test detail columns (E,F) for the
# "change light color" story, widen those columns, tag the header cells as
# wrapped text, and move the active selection onto the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns E (Notes) & F (Test): widen to match the authored layout.
# ColumnWidth is expressed in (rounded) character units on top of a pixel
# grid, so we dial in the input that lands the stored width on the target.
$ws.Columns.Item(5).ColumnWidth = 32.3
$ws.Columns.Item(6).ColumnWidth = 27.17

# --- Header cells E1/F1 pick up word-wrap (bold carries over from the row).
$ws.Range("E1").WrapText = $true
$ws.Range("F1").WrapText = $true

# --- Row 2 grows to fit the new wrapped detail text.
$ws.Rows.Item(2).RowHeight = 45

# --- Fill in the Test / Notes detail for "User can change color of light".
# Write F2 before E2 so the new shared-string entries land in the same
# order as the source workbook (Test string first, then Notes string).
$ws.Range("F2").Value = "1. Can change the color of a light 2. Invalid color input has no effect on color of light"
$ws.Range("E2").Value = "Need to be able to create light, set color, and see result"
$ws.Range("E2").WrapText = $true
$ws.Range("F2").WrapText = $true

# --- Move the active selection to the newly populated area.
$ws.Range("E6").Select()
